$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename "Exam" class header to "Assessment"
$ws.Range("A1").Value = "Assessment"

# Rename field "ExamDescription" -> "AssessmentDescription"
$ws.Range("A7").Value = "AssessmentDescription"

# Rename field "TimeCreated" -> "DateCreated"
$ws.Range("A3").Value = "DateCreated"

# Update the active selection to A4, matching the authored workbook view
$ws.Range("A4").Select()
